$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 241.66667
$ws.Range("J33").Value = 450
$ws.Range("L33").Value = 450
$ws.Range("N33").Value = -908
$ws.Range("H51").Value = 2991.4893
$ws.Range("I51").Value = 2995.652
$ws.Range("J51").Value = 2800
$ws.Range("K51").Value = 2995.652
$ws.Range("L51").Value = 2800
$ws.Range("M51").Value = -2511.652
$ws.Range("N51").Value = -3768
$ws.Range("H137").Value = 31258864
$ws.Range("I137").Value = 38464372
$ws.Range("J137").Value = 34998.332
$ws.Range("K137").Value = 115393116
$ws.Range("L137").Value = 104994.996
$ws.Range("M137").Value = -115390566
$ws.Range("N137").Value = -110094.996
$ws.Range("H138").Value = 11454.556
$ws.Range("I138").Value = 13818.8
$ws.Range("K138").Value = 41456.39999999999
$ws.Range("M138").Value = -36316.39999999999
$ws.Range("H141").Value = 1111.88
$ws.Range("I141").Value = 1111.88
$ws.Range("K141").Value = 3335.64
$ws.Range("M141").Value = 1844.36

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 760.1
$ws.Range("I32").Value = 760.1
$ws.Range("K32").Value = 760.1
$ws.Range("M32").Value = -473.1
$ws.Range("H132").Value = 3731.5667
$ws.Range("I132").Value = 3214.3635
$ws.Range("K132").Value = 9643.0905
$ws.Range("M132").Value = -7113.0905
$ws.Range("H134").Value = 89999
$ws.Range("J134").Value = 89999
$ws.Range("L134").Value = 89999
$ws.Range("N134").Value = -100139

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2425.7
$ws.Range("I86").Value = 2529.6667
$ws.Range("K86").Value = 2529.6667
$ws.Range("M86").Value = -1406.6667
$ws.Range("H89").Value = 2425.7
$ws.Range("I89").Value = 2529.6667
$ws.Range("K89").Value = 12648.3335
$ws.Range("M89").Value = -7032.333500000001
$ws.Range("H99").Value = 10954.272
$ws.Range("I99").Value = 12722
$ws.Range("J99").Value = 2999.5
$ws.Range("K99").Value = 12722
$ws.Range("L99").Value = 2999.5
$ws.Range("M99").Value = -11224
$ws.Range("N99").Value = -5995.5
$ws.Range("H132").Value = 80279.336
$ws.Range("J132").Value = 80279.336
$ws.Range("L132").Value = 80279.336
$ws.Range("N132").Value = -90399.336
$ws.Range("H134").Value = 4391965
$ws.Range("I134").Value = 3332.2424
$ws.Range("K134").Value = 9996.727200000001
$ws.Range("M134").Value = -7461.727200000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()
$ws.Range("H60").Value = 25000
$ws.Range("J60").Value = 31000
$ws.Range("L60").Value = 31000
$ws.Range("N60").Value = -32022
$ws.Range("H62").Value = 5397.5
$ws.Range("I62").Value = 5395
$ws.Range("K62").Value = 5395
$ws.Range("M62").Value = -4771
$ws.Range("H65").Value = 5397.5
$ws.Range("I65").Value = 5395
$ws.Range("K65").Value = 26975
$ws.Range("M65").Value = -23855
$ws.Range("H68").Value = 45000
$ws.Range("J68").Value = 45000
$ws.Range("L68").Value = 45000
$ws.Range("N68").Value = -46498
$ws.Range("H71").Value = 45000
$ws.Range("J71").Value = 45000
$ws.Range("L71").Value = 135000
$ws.Range("N71").Value = -142488
$ws.Range("H99").Value = 53625
$ws.Range("I99").Value = 102250
$ws.Range("K99").Value = 102250
$ws.Range("M99").Value = -100752
$ws.Range("H126").Value = 53625
$ws.Range("I126").Value = 102250
$ws.Range("K126").Value = 306750
$ws.Range("M126").Value = -304280
$ws.Range("H132").Value = 5305.7393
$ws.Range("I132").Value = 5092.077
$ws.Range("J132").Value = 5583.5
$ws.Range("K132").Value = 15276.231
$ws.Range("L132").Value = 16750.5
$ws.Range("M132").Value = -12746.231
$ws.Range("N132").Value = -21810.5
$ws.Range("H134").Value = 1421.6666
$ws.Range("I134").Value = 1341.0123
$ws.Range("J134").Value = 3599.3333
$ws.Range("K134").Value = 4023.0369
$ws.Range("L134").Value = 10797.9999
$ws.Range("M134").Value = -1488.0369
$ws.Range("N134").Value = -15867.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 50.555557
$ws.Range("J2").Value = 39.4
$ws.Range("L2").Value = 236.4
$ws.Range("N2").Value = -462.4
$ws.Range("H8").Value = 417.25
$ws.Range("I8").Value = 417.25
$ws.Range("K8").Value = 1251.75
$ws.Range("M8").Value = -1112.75
$ws.Range("H37").Value = 77348.95
$ws.Range("J37").Value = 77348.95
$ws.Range("L37").Value = 232046.85
$ws.Range("N37").Value = -232270.85
$ws.Range("H68").Value = 2265.2
$ws.Range("I68").Value = 1798.125
$ws.Range("J68").Value = 2485
$ws.Range("K68").Value = 5394.375
$ws.Range("L68").Value = 7455
$ws.Range("M68").Value = -4583.375
$ws.Range("N68").Value = -9077
$ws.Range("H71").Value = 2265.2
$ws.Range("I71").Value = 1798.125
$ws.Range("J71").Value = 2485
$ws.Range("K71").Value = 16183.125
$ws.Range("L71").Value = 22365
$ws.Range("M71").Value = -12127.125
$ws.Range("N71").Value = -30477

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 112
$ws.Range("I31").Value = 112
$ws.Range("K31").Value = 112
$ws.Range("M31").Value = 180
$ws.Range("H37").Value = 112
$ws.Range("I37").Value = 112
$ws.Range("K37").Value = 112
$ws.Range("M37").Value = 165
$ws.Range("H132").Value = 25575
$ws.Range("I132").Value = 26985.062
$ws.Range("K132").Value = 80955.186
$ws.Range("M132").Value = -78425.186

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H42").Value = 18955.2
$ws.Range("J42").Value = 18955.2
$ws.Range("L42").Value = 18955.2
$ws.Range("N42").Value = -20081.2
$ws.Range("H49").Value = 18955.2
$ws.Range("J49").Value = 18955.2
$ws.Range("L49").Value = 18955.2
$ws.Range("N49").Value = -19249.2
$ws.Range("H68").Value = 12146.3
$ws.Range("I68").Value = 10325.632
$ws.Range("K68").Value = 10325.632
$ws.Range("M68").Value = -9576.632
$ws.Range("H71").Value = 12146.3
$ws.Range("I71").Value = 10325.632
$ws.Range("K71").Value = 51628.16
$ws.Range("M71").Value = -47884.16
$ws.Range("H132").Value = 1628451.6
$ws.Range("I132").Value = 2022149.1
$ws.Range("J132").Value = 4449.75
$ws.Range("K132").Value = 6066447.300000001
$ws.Range("L132").Value = 13349.25
$ws.Range("M132").Value = -6063917.300000001
$ws.Range("N132").Value = -18409.25
$ws.Range("H133").Value = 89298.664
$ws.Range("J133").Value = 89298.664
$ws.Range("L133").Value = 89298.664
$ws.Range("N133").Value = -94358.664

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 3000000
$ws.Range("I8").Value = 3000000
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 3000000
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -2999860
$ws.Range("N8").ClearContents()
$ws.Range("H41").Value = 32911.355
$ws.Range("J41").Value = 33135.31
$ws.Range("L41").Value = 33135.31
$ws.Range("N41").Value = -33915.31
$ws.Range("H54").Value = 26000
$ws.Range("I54").Value = 25000
$ws.Range("J54").Value = 28000
$ws.Range("K54").Value = 25000
$ws.Range("L54").Value = 28000
$ws.Range("M54").Value = -24480
$ws.Range("N54").Value = -29040
$ws.Range("H132").Value = 5209325.5
$ws.Range("I132").Value = 6173845
$ws.Range("J132").Value = 920
$ws.Range("K132").Value = 18521535
$ws.Range("L132").Value = 2760
$ws.Range("M132").Value = -18519005
$ws.Range("N132").Value = -7820
